$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, shifting existing rows 88..220 down to 89..221
$ws.Rows("88:88").Insert()

# Populate the newly inserted row 88 with the new record's data
$ws.Range("A88").Value = 5
$ws.Range("B88").Value = "Macroferia Regional de Talca"
$ws.Range("C88").Value = "Maule"
$ws.Range("D88").Value = 44557
$ws.Range("E88").Value = 7
$ws.Range("F88").Value = 100112003
$ws.Range("G88").Value = "Ajo"
$ws.Range("H88").Value = "Chino"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 300
$ws.Range("K88").Value = 18000
$ws.Range("L88").Value = 18000
$ws.Range("M88").Value = 18000
$ws.Range("N88").Value = "$/caja 10 kilos"
$ws.Range("O88").Value = "China"
$ws.Range("P88").Value = 1800
$ws.Range("Q88").Value = 10
$ws.Range("R88").Value = "Hortaliza"
